# Fill in the camera-model (E) / battery-model (F) lookup text that was
# typed into rows 29-37 of the "Good Looks" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Good Looks")

$ws.Range("E29").Value = "Fujifilm X-A1 Fujifilm X-A2 Fujifilm X-E1 Fujifilm X-E2 Fujifilm X-M1 Fujifilm X-Pro1 Fujifilm X-T1 Fujifilm X-T1 IR Fujifilm X-T10 Fujifilm FinePix HS30EXR Fujifilm FinePix HS33EXR Fujifilm FinePix HS35EXR Fujifilm FinePix HS50EXR"
$ws.Range("F29").Value = "Fuji NP-W126, NPW126"

$ws.Range("E30").Value = "Samsung EX2F Samsung L100, L110, L120, L200, L210, L310W, L313 Samsung ES50, ES55 Samsung HZ10W, HZ15W Samsung IT100 Samsung M100, M110, M310W Samsung NV9 Samsung P800, P1000 Samsung PL50, PL51, PL55, PL57, PL60, PL65, PL70 Samsung SL102, SL105, SL202, SL203, SL310, SL310W, SL420, SL502, SL620, SL720, SL820 Samsung TL9 Samsung WB150F Samsung WB250F Samsung WB350F Samsung WB380 Samsung WB500 Samsung WB550 Samsung WB750 Samsung WB850F Samsung WB800F Samsung WB1100F Samsung WB1101F Samsung WB1102 Samsung WB2100 Samsung HMX-U100 Samsung HMX-U100UN, HMXU100UN, U100UN Samsung HMX-U100SN, HMXU100SN, U100SN"
$ws.Range("F30").Value = "Samsung SLB-10A, SLB10A, 4302-001221, 4302001221"

$ws.Range("E31").Value = "AG-AC7, AGAC7, AC7 AG-HMC40, AGHMC40, HMC40 AG-HMC70, AGHMC70, HMC70 AG-HMC150, AGHMC150, HMC150 HDC-DX1, HDCDX1, HDCDX1 HDC-DX3, HDCDX3, HDCDX3 HDC-HS100 HDC-HS100GK, HDCHS100GK, HS100GK HDC-HS20 HDC-HS20K, HDCHS20K, HS20K HDC-HS200 HDC-HS200K, HDCHS200K, HS200K HDC-HS250 HDC-HS250K, HDCHS250K, HS250K HDC-HS300 HDC-HS300K, HDCHS300K, HS300K HDC-HS700 HDC-HS700K, HDCHS700K, HS700K HDC-HS9, HDCHS9, HS9 HDC-SD100 HDC-SD100GK, HDCSD100GK, SD100GK HDC-SD200, HDCSD200, SD200 HDC-SD20 HDC-SD20K, HDCSD20K, SD20K HDC-SD3, HDCSD3, SD3 HDC-SD9, HDCSD9, SD9 HDC-SD600 HDC-SD600K, HDCSD600K, SD600K HDC-SD700 HDC-SD700K, HDCSD700K, SD700K HDC-SX5, HDCSX5, SX5 HDC-SDT750 HDC-SDT750K, HDCSDT750K, SDT750K HDC-TM20 HDC-TM20K, HDCTM20K, TM20K HDC-TM20K8, HDCTM20K8, TM20K8 HDC-TM20R, HDCTM20R, TM20R HDC-TM20S, HDCTM20S, TM20S HDC-TM300 HDC-TM300K, HDCTM300K, TM300K HDC-TM700 HDC-TM700K, HDCTM700K, TM700K PV-GS90, PVGS90, GS90 SDR-H40 SDR-H41 SDR-H50 SDR-H60 SDR-H79 SDR-H79K, SDRH79K, H79K SDR-H79P, SDRH79P, H79P SDR-H80 SDR-H80S, SDRH80S, H80S SDR-H80K, SDRH80K, H80K SDR-H80A, SDRH80A, H80A SDR-H80R, SDRH80R, H80R SDR-H90 SDR-H90K, SDRH90K, H90K VDR-D50 VDR-D50P, VDRD50P, D50P VDR-D51, VDRD51, D51"
$ws.Range("F31").Value = "VW-VBG260, VW-VBG260K, VW-VBG260E, VW-VBG260PP, VWVBG260, VWVBG260K, VWVBG260E, VWVBG260PP"

$ws.Range("E32").Value = "AG-AC7, AGAC7, AC7 AG-AC130A, AGAC130A, AC130A AG-AC160A, AGAC160A, AC160A AG-HMC40, AGHMC40, HMC40 AG-HMC70, AGHMC70, HMC70 AG-HMC150, AGHMC150, HMC150 HDC-DX1, HDCDX1, HDCDX1 HDC-DX3, HDCDX3, HDCDX3 HDC-HS100 HDC-HS100GK, HDCHS100GK, HS100GK HDC-HS20 HDC-HS20K, HDCHS20K, HS20K HDC-HS200 HDC-HS200K, HDCHS200K, HS200K HDC-HS250 HDC-HS250K, HDCHS250K, HS250K HDC-HS300 HDC-HS300K, HDCHS300K, HS300K HDC-HS700 HDC-HS700K, HDCHS700K, HS700K HDC-HS9, HDCHS9, HS9 HDC-SD100 HDC-SD100GK, HDCSD100GK, SD100GK HDC-SD200, HDCSD200, SD200 HDC-SD20 HDC-SD20K, HDCSD20K, SD20K HDC-SD3, HDCSD3, SD3 HDC-SD9, HDCSD9, SD9 HDC-SD600 HDC-SD600K, HDCSD600K, SD600K HDC-SD700 HDC-SD700K, HDCSD700K, SD700K HDC-SX5, HDCSX5, SX5 HDC-SDT750 HDC-SDT750K, HDCSDT750K, SDT750K HDC-TM20 HDC-TM20K, HDCTM20K, TM20K HDC-TM20K8, HDCTM20K8, TM20K8 HDC-TM20R, HDCTM20R, TM20R HDC-TM20S, HDCTM20S, TM20S HDC-TM300 HDC-TM300K, HDCTM300K, TM300K HDC-TM700 HDC-TM700K, HDCTM700K, TM700K PV-GS90, PVGS90, GS90 SDR-H79 SDR-H79K, SDRH79K, H79K SDR-H79P, SDRH79P, H79P SDR-H80 SDR-H80S, SDRH80S, H80S SDR-H80K, SDRH80K, H80K SDR-H80A, SDRH80A, H80A SDR-H80R, SDRH80R, H80R SDR-H90 SDR-H90K, SDRH90K, H90K VDR-D50 VDR-D50P, VDRD50P, D50P VDR-D51, VDRD51, D51"
$ws.Range("F32").Value = "VW-VBG6, VW-VBG6K, VW-VBG6E, VW-VBG6PP, VWVBG6, VWVBG6K, VWVBG6E, VWVBG6PP"

$ws.Range("E33").Value = "HC-V10 HC-V10K, HCV10K, V10K HC-V100 HC-V100K, HCV100K, V100K HC-V100M, HCV100M, V100M HC-V500 HC-V500K, HCV500K, V500K HC-V500M, HCV500M, V500M HDC-SD40 HDC-SD40K, HDCSD40K, SD40K HDC-SD60 HDC-SD60K, HDCSD60K, SD60K HDC-SD60S, HDCSD60S, SD60S HDC-SD80 HDC-SD80K, HDCSD80K, SD80K HDC-SD80S, HDCSD80S, SD80S HDC-SD80R, HDCSD80R, SD80R HDC-SD90 HDC-SD90K, HDCSD90K, SD90K HDC-SDX1 HDC-SDX1H, HDCSDX1H, SDX1H HDC-TMX1, HDCTMX1, TMX1 HDC-TM40 HDC-TM40K, HDCTM40K, TM40K HDC-TM41 HDC-TM41H, HDCTM41H, TM41H HDC-TM55 HDC-TM55K, HDCTM55K, TM55K HDC-TM60 HDC-TM60K, HDCTM60K, TM60K HDC-TM80 HDC-TM80K, HDCTM80K, TM80K HDC-TM80S, HDCTM80S, TM80S HDC-TM80R, HDCTM80R, TM80R HDC-TM90 HDC-TM90K, HDCTM90K, TM90K HDC-TMX1, HDCTMX1, TMX1 HDC-HS40 HDC-HS40K, HDCHS40K, HS40K HDC-HS60 HDC-HS60K, HDCHS60K, HS60K HDC-HS80 HDC-HS80K, HDCHS80K, HS80K SDR-H85 SDR-H85A, SDRH85A, H85A SDR-H85K, SDRH85K, H85K SDR-H85S, SDRH85S, H85S SDR-H95 SDR-H95K, SDRH95K, H95K SDR-H100 SDR-H100K, SDRH100K, H100K SDR-H100S, SDRH100S, H100S SDR-H100R, SDRH100R, H100R SDR-H101 SDR-H101K, SDRH101K, H101K SDR-H101R, SDRH101R, H101R SDR-HS60 SDR-HS60K, SDRHS60K, HS60K SDR-HS80 SDR-HS80K, SDRHS80K, HS80K SDR-S45, SDRS45, S45 SDR-S50 SDR-S50A, SDRS50A, S50A SDR-S50K, SDRS50K, S50K SDR-S50N, SDRS50N, S50N SDR-S70 SDR-S70K, SDRS70K, S70K SDR-S70S, SDRS70S, S70S SDR-S70R, SDRS70R, S70R SDR-S71 SDR-S71K, SDRS71K, S71K SDR-T50 SDR-T50K, SDRT50K, T50K SDR-T55 SDR-T55K, SDRT55K, T55K SDR-T70 SDR-T70K, SDRT70K, T70K SDR-T71 SDR-T71K, SDRT71K, T71K SDR-T76 SDR-T76K, SDRT76K, T76K"
$ws.Range("F33").Value = "VW-VBL090, VWVBL090, VBL090, VW-VBL90, VWVBL90, VBL90, VW-VBK180, VW-VBK180K, VW-VBK180PP, VWVBK180, VWVBK180K, VWVBK180PP, VW-VBK180E, VWVBK180E, VW-VBK360, VW-VBK360K, VW-VBK360PP, VWVBK360, VWVBK360K, VWVBK360PP, VW-VBK360E, VWVBK360E"

$ws.Range("E34").Value = $ws.Range("E33").Text
$ws.Range("F34").Value = $ws.Range("F33").Text

$ws.Range("E35").Value = "HC-X900, HCX900, X900 HC-X920, HCX920, X920 HDC-HS900 HDC-HS900K, HDCHS900K, HS900K HDC-SD800 HDC-SD800K, HDCSD800K, SD800K HDC-SD900 HDC-SD900K, HDCSD900K, SD900K HDC-SD909 HDC-SD909K, HDCSD909K, SD909K HDC-TM900 HDC-TM900K, HDCTM900K, TM900K"
$ws.Range("F35").Value = "VW-VBN130, VW-VBN130K, VW-VBN130PP, VW-VBN130E, VWVBN130, VWVBN130K, VWVBN130PP, VWVBN130E"

$ws.Range("E36").Value = $ws.Range("E35").Text
$ws.Range("F36").Value = "VW-VBN260, VW-VBN260K, VW-VBN260PP, VW-VBN260E, VWVBN260, VWVBN260K, VWVBN260PP, VWVBN260E"

$ws.Range("E37").Value = "GC-PX10, GCPX10, PX10 GC-PX100, GCPX100, PX100 GR-D720US, GRD720US, D720US GR-D720EK, GRD720EK, D720EK GR-D720EX, GRD720EX, D720EX GR-D720, GRD720 GR-D725US, GRD725US, D725US GR-D725EK, GRD725EK, D725EK GR-D725EX, GRD725EX, D725EX GR-D725, GRD725 GR-D726US, GRD726US, D726US GR-D726EK, GRD726EK, D726EK GR-D726EX, GRD726EX, D726EX GR-D726, GRD726 GR-D728US, GRD728US, D728US GR-D728EK, GRD728EK, D728EK GR-D728EX, GRD728EX, D728EX GR-D728, GRD728 GR-D740US, GRD740US, D740US GR-D740EK, GRD740EK, D740EK GR-D740EX, GRD740EX, D740EX GR-D740, GRD740 GR-D745US, GRD745US, D745US GR-D745EK, GRD745EK, D745EK GR-D745EX, GRD745EX, D745EX GR-D745, GRD745 GR-D746US, GRD746US, D746US GR-D746EK, GRD746EK, D746EK GR-D746EX, GRD746EX, D746EX GR-D746, GRD746 GR-D750US, GRD750US, D750US GR-D750EK, GRD750EK, D750EK GR-D750EX, GRD750EX, D750EX GR-D750, GRD750 GR-D760US, GRD760US, D760US GR-D760EK, GRD760EK, D760EK GR-D760EX, GRD760EX, D760EX GR-D760, GRD760 GR-D770US, GRD770US, D770US GR-D770EK, GRD770EK, D770EK GR-D770EX, GRD770EX, D770EX GR-D770, GRD770 GR-D770E, GRD770E, D770E GR-D770VS, GRD770VS, D770VS GR-D771US, GRD771US, D771US GR-D771EK, GRD771EK, D771EK GR-D771EX, GRD771EX, D771EX GR-D771, GRD771 GR-D775US, GRD775US, D775US GR-D775EK, GRD775EK, D775EK GR-D775EX, GRD775EX, D775EX GR-D775, GRD775 GR-D790US, GRD790US, D790US GR-D790EK, GRD790EK, D790EK GR-D790EX, GRD790EX, D790EX GR-D790, GRD790 GR-D796US, GRD796US, D796US GR-D796EK, GRD796EK, D796EK GR-D796EX, GRD796EX, D796EX GR-D796, GRD796 GR-D850US, GRD850US, D850US GR-D850EK, GRD850EK, D850EK GR-D850EX, GRD850EX, D850EX GR-D850, GRD850 GR-D870US, GRD870US, D870US GR-D870EK, GRD870EK, D870EK GR-D870EX, GRD870EX, D870EX GR-D870, GRD870 GR-D875US, GRD875US, D875US GR-D875EK, GRD875EK, D875EK GR-D875EX, GRD875EX, D875EX GR-D875, GRD875 GR-DA30US, GRDA30US, DA30US GR-DA30, GRDA30 GS-TD1, GSTD1, TD1 GY-HM70U, GYHM70U GY-HM100, GYHM100 GY-HM100U, GYHM100U GY-HM100E, GYHM100E GY-HM100EC, GYHM100EC GY-HM150, GYHM150 GY-HM150U, GYHM150U GY-HM150E, GYHM150E GY-HMZ1, GYHMZ1 GY-HMZ1E, GYHMZ1E GY-HMZ1U, GYHMZ1U Everio GZ-MG130US, GZMG130US, MG130US Everio GZ-MG130EK, GZMG130EK, MG130EK Everio GZ-MG130EX, GZMG130EX, MG130EX Everio GZ-MG130, GZMG130 Everio GZ-MG131US, GZMG131US, MG131US Everio GZ-MG131EK, GZMG131EK, MG131EK Everio GZ-MG131EX, GZMG131EX, MG131EX Everio GZ-MG131, GZMG131 Everio GZ-MG132US, GZMG132US, MG132US Everio GZ-MG132EK, GZMG132EK, MG132EK Everio GZ-MG132EX, GZMG132EX, MG132EX Everio GZ-MG132, GZMG132 Everio GZ-MG133US, GZMG133US, MG133US Everio GZ-MG133EK, GZMG133EK, MG133EK Everio GZ-MG133EX, GZMG133EX, MG133EX Everio GZ-MG133, GZMG133 Everio GZ-MG134US, GZMG134US, MG134US Everio GZ-MG134EK, GZMG134EK, MG134EK Everio GZ-MG134EX, GZMG134EX, MG134EX Everio GZ-MG134, GZMG134 Everio GZ-MG135US, GZMG135US, MG135US Everio GZ-MG135EK, GZMG135EK, MG135EK Everio GZ-MG135EX, GZMG135EX, MG135EX Everio GZ-MG135, GZMG135 Everio GZ-MG148US, GZMG148US, MG148US Everio GZ-MG148EK, GZMG148EK, MG148EK Everio GZ-MG148EX, GZMG148EX, MG148EX Everio GZ-MG148, GZMG148 Everio GZ-MG150US, GZMG150US, MG150US Everio GZ-MG150EK, GZMG150EK, MG150EK Everio GZ-MG150EX, GZMG150EX, MG150EX Everio GZ-MG150, GZMG150 Everio GZ-MG155US, GZMG155US, MG155US Everio GZ-MG155EK, GZMG155EK, MG155EK Everio GZ-MG155EX, GZMG155EX, MG155EX Everio GZ-MG155, GZMG155 Everio GZ-MG157US, GZMG157US, MG157US Everio GZ-MG157EK, GZMG157EK, MG157EK Everio GZ-MG157EX, GZMG157EX, MG157EX Everio GZ-MG157, GZMG157 Everio GZ-MG175US, GZMG175US, MG175US Everio GZ-MG175EK, GZMG175EK, MG175EK Everio GZ-MG175EX, GZMG175EX, MG175EX Everio GZ-MG175, GZMG175 Everio GZ-MG177US, GZMG177US, MG177US Everio GZ-MG177EK, GZMG177EK, MG177EK Everio GZ-MG177EX, GZMG177EX, MG177EX Everio GZ-MG177, GZMG177 Everio GZ-MG230, GZMG230, MG230 Everio GZ-MG255US, GZMG255US, MG255US Everio GZ-MG255EK, GZMG255EK, MG255EK Everio GZ-MG255EX, GZMG255EX, MG255EX Everio GZ-MG255, GZMG255 Everio GZ-MG275US, GZMG275US, MG275US Everio GZ-MG275EK, GZMG275EK, MG275EK Everio GZ-MG275EX, GZMG275EX, MG275EX Everio GZ-MG275, GZMG275 Everio GZ-MG275E, GZMG275E Everio GZ-MG330, GZMG330 Everio GZ-MG330A, GZMG330A Everio GZ-MG330H, GZMG330H Everio GZ-MG330R, GZMG330R Everio GZ-MG330HUS, GZMG330HUS, MG330HUS Everio GZ-MG330RUS, GZMG330RUS, MG330RUS Everio GZ-MG330AUS, GZMG330AUS, MG330US Everio GZ-MG335, GZMG335 Everio GZ-MG335H, GZMG335H Everio GZ-MG335W, GZMG335W Everio GZ-MG335HUS, GZMG335HUS, MG335HUS Everio GZ-MG335WUS, GZMG335WUS, MG335WUS Everio GZ-MG340, GZMG340 Everio GZ-MG340B, GZMG340B Everio GZ-MG340BUS, GZMG340BUS, MG340BUS Everio GZ-MG360, GZMG360 Everio GZ-MG360B, GZMG360B Everio GZ-MG360BUS, GZMG360BUS, MG360BUS Everio GZ-MG365, GZMG365 Everio GZ-MG365B, GZMG365B Everio GZ-MG365BUS, GZMG365BUS, MG365BUS Everio GZ-MG430, GZMG430 Everio GZ-MG430H, GZMG430H Everio GZ-MG430HUS, GZMG430HUS, MG430HUS Everio GZ-MG435, GZMG435 Everio GZ-MG435H, GZMG435H Everio GZ-MG435HUS, GZMG435HUS, MG435HUS Everio GZ-MG465, GZMG465 Everio GZ-MG465B, GZMG465B Everio GZ-MG465BUS, GZMG465BUS, MG465BUS Everio GZ-MG555US, GZMG555US, MG555US Everio GZ-MG555EK, GZMG555EK, MG555EK Everio GZ-MG555EX, GZMG555EX, MG555EX Everio GZ-MG555, GZMG555 Everio GZ-MG575US, GZMG575US, MG575US Everio GZ-MG575EK, GZMG575EK, MG575EK Everio GZ-MG575EX, GZMG575EX, MG575EX Everio GZ-MG575, GZMG575 GZ-MG630A, GZMG630A, MG630A GZ-MG630R, GZMG630R, MG630R GZ-MG630S, GZMG630S, MG630S GZ-MG630US, GZMG630US, MG630US GZ-MG670B, GZMG670B, MG670B GZ-MG670US, GZMG670US, MG670US GZ-MG680B, GZMG680B, MG680B GZ-MG680US, GZMG680US, MG680US GZ-MG730US, GZMG730US, MG730US Everio GZ-MS100, GZMS100 Everio GZ-MS100US, GZMS100US, MS100US Everio GZ-MS100RUS, GZMS100RUS, MS100RUS Everio GZ-MS100R, GZMS100R Everio GZ-MS120, GZMS120, MS120 Everio GZ-MS120A, GZMS120A, MS120A Everio GZ-MS120R, GZMS120R, MS120R Everio GZ-MS120B, GZMS120B, MS120B Everio GZ-MS130, GZMS130 Everio GZ-MS130A, GZMS130A, MS130A Everio GZ-MS130R, GZMS130R, MS130R Everio GZ-MS130B, GZMS130B, MS130B Everio GZ-HD3US, GZHD3US, HD3US Everio GZ-HD3EK, GZHD3EK, HD3EK Everio GZ-HD3EX, GZHD3EX, HD3EX Everio GZ-HD3, GZHD3 Everio GZ-HD5US, GZHD5US, HD5US Everio GZ-HD5EK, GZHD5EK, HD5EK Everio GZ-HD5EX, GZHD5EX, HD5EX Everio GZ-HD5, GZHD5 Everio GZ-HD6US, GZHD6US, HD6US Everio GZ-HD6EK, GZHD6EK, HD6EK Everio GZ-HD6EX, GZHD6EX, HD6EX Everio GZ-HD6, GZHD6 Everio GZ-HD7US, GZHD7US, HD7US Everio GZ-HD7EK, GZHD7EK, HD7EK Everio GZ-HD7EX, GZHD7EX, HD7EX Everio GZ-HD7, GZHD7 Everio GZ-HD10US, GZHD10US, HD10US Everio GZ-HD10EK, GZHD10EK, HD10EK Everio GZ-HD10EX, GZHD10EX, HD10EX Everio GZ-HD10, GZHD10 Everio GZ-HD30US, GZHD30US, HD30US Everio GZ-HD30EK, GZHD30EK, HD30EK Everio GZ-HD30EX, GZHD30EX, HD30EX Everio GZ-HD30, GZHD30 Everio GZ-HD40US, GZHD40US, HD40US Everio GZ-HD40EK, GZHD40EK, HD40EK Everio GZ-HD40EX, GZHD40EX, HD40EX Everio GZ-HD40, GZHD40 GZ-HD300, GZHD300 GZ-HD300B, GZHD300B, HD300B GZ-HD300R, GZHD300R, HD300R GZ-HD300A, GZHD300A, HD300A GZ-HD200B, GZHD200B, HD200B GZ-HD200R, GZHD200R, HD200R GZ-HD200A, GZHD200A, HD200A GZ-HM1, GZHM1, HM1 GZ-HM200, GZHM200, HM200 GZ-HM400, GZHM400, HM400 GZ-X900, GZX900, X900 Victor JVC BN-VF815U, BN-VF815, BNVF815U, BNVF815, BN-VF815US, BNVF815US, BN-VF815USM, BNVF815USM, BN-VF815JP, BNVF815JP"

# Match the author's last cursor position / selection captured in the diff
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("E38").Select()
